# Apply line-break splits to four paragraphs per the target diff.
$d = $word.ActiveDocument

# 1) "Programa" section - split into 7 numbered items
$d.Content.Find.Execute(
    '1.Introdução ao projeto em Engenharia: o que é projeto em engenharia e por que projetar? Metodologias de projeto; etapas de elaboração de projeto;2.Metodologia de projeto focada no ser humano Design Thinking. Entendimento do duplo diamante da inovação. Etapas do Design Thinking: empatia, definição do problema, ideação, prototipação do plano e teste do produto;3. Processo de melhoria contínua Kaizen. Ciclo de vida de projeto PDCA (Plan-Do-Check-Act): Planejar-Desenvolver-Checar-Agir;4.Métodos e normas para redação de textos científicos;5.Desenvolvimento de um projeto temático, compreendendo: definição do problema e formação de alternativas de solução; estabelecimento de critérios; escolha e avaliação de soluções; especificação da solução;6.Noções de planejamento e gestão de projetos; organização do tempo; técnicas para a realização de apresentações; noções de aprendizagem baseada em projetos; trabalho em grupo, equipes e times7.Tutoria de projetos',
    $false, $false, $false, $false, $false, $true, 1, $false,
    '1.Introdução ao projeto em Engenharia: o que é projeto em engenharia e por que projetar? Metodologias de projeto; etapas de elaboração de projeto;^l2.Metodologia de projeto focada no ser humano Design Thinking. Entendimento do duplo diamante da inovação. Etapas do Design Thinking: empatia, definição do problema, ideação, prototipação do plano e teste do produto;^l3. Processo de melhoria contínua Kaizen. Ciclo de vida de projeto PDCA (Plan-Do-Check-Act): Planejar-Desenvolver-Checar-Agir;^l4.Métodos e normas para redação de textos científicos;^l5.Desenvolvimento de um projeto temático, compreendendo: definição do problema e formação de alternativas de solução; estabelecimento de critérios; escolha e avaliação de soluções; especificação da solução;^l6.Noções de planejamento e gestão de projetos; organização do tempo; técnicas para a realização de apresentações; noções de aprendizagem baseada em projetos; trabalho em grupo, equipes e times^l7.Tutoria de projetos',
    2
)

# 2) "Método:" run text - split into 4 sentences
$d.Content.Find.Execute(
    'O método utilizado tem por fundamento a aprendizagem baseada em projetos que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, capacidade de comunicação oral e verbal e aspectos de liderança, dentre outros.Os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Materiais, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão.Cada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.As aulas ocorrerão por meio de uma reunião da equipe de trabalho para tratar do projeto; palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores ou profissionais de empresas.',
    $false, $false, $false, $false, $false, $true, 1, $false,
    'O método utilizado tem por fundamento a aprendizagem baseada em projetos que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, capacidade de comunicação oral e verbal e aspectos de liderança, dentre outros.^lOs alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Materiais, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão.^lCada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.^lAs aulas ocorrerão por meio de uma reunião da equipe de trabalho para tratar do projeto; palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores ou profissionais de empresas.',
    2
)

# 3) "Critério:" run text - split into 2 sentences
$d.Content.Find.Execute(
    'A nota será individual e será a média ponderada de entregas do projeto, tais como: projeto preliminar, projeto final, envolvimento do aluno com o projeto, avaliação dos pares, autoavaliação e apresentação de trabalhos, dentre outros.O detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na avaliação da disciplina.',
    $false, $false, $false, $false, $false, $true, 1, $false,
    'A nota será individual e será a média ponderada de entregas do projeto, tais como: projeto preliminar, projeto final, envolvimento do aluno com o projeto, avaliação dos pares, autoavaliação e apresentação de trabalhos, dentre outros.^lO detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na avaliação da disciplina.',
    2
)

# 4) Bibliography paragraph - split into 5 references
$d.Content.Find.Execute(
    '- BAZZO, Walter; PEREIRA, Luiz T.V. Introdução à Engenharia, 3a. edição. Florianópolis: Editora da UFSC, 2013.- COCIAN, Luis Fernando Espinosa. Introdução à Engenharia. Porto Alegre: Bookman, 2017.- BENNETT, Ronald; MILLAM, Elaine. Liderança para engenheiros. Porto Alegre: AMGH, 2014.- ALEXANDER, C. K.; WATSON, J. A. Habilidades para uma carreira de sucesso na engenharia, Porto Alegre: AMGH Editora, 2015.- MCCAHAN, S.; ANDERSON, P.; KORTSCHOT, M.; WEISS, P.; WOODHOUSE, K. Projetos de Engenharia: uma introdução. 1ª edição. -Rio de Janeiro: LTC, 2017.',
    $false, $false, $false, $false, $false, $true, 1, $false,
    '- BAZZO, Walter; PEREIRA, Luiz T.V. Introdução à Engenharia, 3a. edição. Florianópolis: Editora da UFSC, 2013.^l- COCIAN, Luis Fernando Espinosa. Introdução à Engenharia. Porto Alegre: Bookman, 2017.^l- BENNETT, Ronald; MILLAM, Elaine. Liderança para engenheiros. Porto Alegre: AMGH, 2014.^l- ALEXANDER, C. K.; WATSON, J. A. Habilidades para uma carreira de sucesso na engenharia, Porto Alegre: AMGH Editora, 2015.^l- MCCAHAN, S.; ANDERSON, P.; KORTSCHOT, M.; WEISS, P.; WOODHOUSE, K. Projetos de Engenharia: uma introdução. 1ª edição. -Rio de Janeiro: LTC, 2017.',
    2
)
